$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# The original Sheet1!A1 was an empty placeholder cell; the published
# workbook no longer carries it.
$ws1.Range("A1").ClearContents()

# Add the new "Terms" sheet right after Sheet1 (becomes the active tab,
# exactly like the published workbook).
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Terms"

$ws2.Range("A1").Value = "This dataset on 'SleepData.xlsx' is hypothetical and was generated"
$ws2.Range("A2").Value = "by Paolo G. Hilado (Github: Dcroix) for training purposes on Basic Statistics . Considering"
$ws2.Range("A3").Value = "that most of the values generated by this dataset use randomization, "
$ws2.Range("A4").Value = "in such a rare case that it resembles any existing dataset, it is purely "
$ws2.Range("A5").Value = "coincidental. It is distributed under "
$ws2.Range("A6").Value = " Creative Commons Attribution-NoDerivatives 4.0 International Public License."

$ws2.Columns.Item(1).ColumnWidth = 70.6666667

$ws2.Hyperlinks.Add($ws2.Range("A6"), "https://creativecommons.org/licenses/by-nd/4.0/")
